$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.187.83"
$ws.Range("E2").Value = '  -2.08%  '

$ws.Range("D3").Value = "'3.684.91"
$ws.Range("E3").Value = '  -2.99%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = "'682.31"
$ws.Range("E5").Value = '  -3.43%  '

$ws.Range("D6").Value = "'162.71"
$ws.Range("E6").Value = '  -4.66%  '

$ws.Range("D7").Value = "'3.684.82"
$ws.Range("E7").Value = '  -2.96%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").Value = '  -4.44%  '

$ws.Range("E10").Value = '  -7.55%  '

$ws.Range("E11").Value = '  -2.16%  '

$ws.Range("E12").Value = '  -1.34%  '

$ws.Range("E13").Value = '  -6.92%  '

$ws.Range("D14").Value = "'33.60"
$ws.Range("E14").Value = '  -7.15%  '

$ws.Range("D15").Value = "'4.303.58"
$ws.Range("E15").Value = '  -3.06%  '

$ws.Range("D16").Value = "'3.678.91"
$ws.Range("E16").Value = '  -2.52%  '

$ws.Range("D17").Value = "'69.246.72"
$ws.Range("E17").Value = '  -2.02%  '

$ws.Range("E18").Value = '  -1.79%  '

$ws.Range("D19").Value = "'16.36"
$ws.Range("E19").Value = '  -5.97%  '

$ws.Range("D20").Value = "'6.66"
$ws.Range("E20").Value = '  -6.81%  '

$ws.Range("D21").Value = "'482.74"
$ws.Range("E21").Value = '  -2.23%  '

$ws.Range("D22").Value = "'9.81"
$ws.Range("E22").Value = '  -7.71%  '

$ws.Range("D23").Value = "'0.667"
$ws.Range("E23").Value = '  -8.52%  '

$ws.Range("D24").Value = "'79.86"
$ws.Range("E24").Value = '  -5.98%  '

$ws.Range("D25").Value = "'3.827.69"
$ws.Range("E25").Value = '  -3.08%  '

$ws.Range("D26").Value = "'11.57"
$ws.Range("E26").Value = '  -4.24%  '

$ws.Range("D27").Value = "'0.0000128"
$ws.Range("E27").Value = '  -12.02%  '

$ws.Range("E28").Value = '  -0.02%  '

$ws.Range("D29").Value = "'9.55"
$ws.Range("E29").Value = '  -8.89%  '

$ws.Range("E30").Value = '  -10.03%  '

$ws.Range("D31").Value = "'2.75"
$ws.Range("E31").Value = '  -11.00%  '

$ws.Range("E32").Value = '  -4.82%  '

$ws.Range("D33").Value = "'6.78"
$ws.Range("E33").Value = '  -7.62%  '

$ws.Range("D34").Value = "'27.00"
$ws.Range("E34").Value = '  -7.28%  '

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  +0.09%  '

$ws.Range("D36").Value = "'0.164"
$ws.Range("E36").Value = '  -6.30%  '

$ws.Range("D37").Value = "'3.650.14"
$ws.Range("E37").Value = '  -3.16%  '

$ws.Range("D38").Value = "'8.53"
$ws.Range("E38").Value = '  -5.91%  '

$ws.Range("D39").Value = "'6.06"
$ws.Range("E39").Value = '  +2.25%  '

$ws.Range("D40").Value = "'0.0946"
$ws.Range("E40").Value = '  -6.63%  '

$ws.Range("E42").Value = '  -6.14%  '

$ws.Range("E43").Value = '  -0.03%  '

$ws.Range("E44").Value = '  -7.95%  '

$ws.Range("D45").Value = "'157.02"
$ws.Range("E45").Value = '  -4.48%  '

$ws.Range("D46").Value = "'48.09"
$ws.Range("E46").Value = '  -1.62%  '

$ws.Range("D47").Value = "'2.79"
$ws.Range("E47").Value = '  -15.23%  '

$ws.Range("D48").Value = "'396.12"
$ws.Range("E48").Value = '  -6.31%  '

$ws.Range("B49").Value = 'FLOKI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D49").Value = "'0.000277"
$ws.Range("E49").Value = '  -13.24%  '

$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").Value = "'1.31"
$ws.Range("E50").Value = '  -4.31%  '

$ws.Range("D51").Value = "'8.09"
$ws.Range("E51").Value = '  -6.90%  '
